# Report.docx edit:
#  1. Insert a new "Benchmark" sub-heading paragraph (underlined) right
#     before the "To benchmark and validate..." paragraph, following the
#     existing "Benchmark and Validation" heading.
#  2. Insert a new "Validation" sub-heading paragraph (underlined) right
#     before the "Furthermore, to validate..." paragraph.
#  3. Insert a new centered, bold, underlined "Results" heading paragraph
#     right before the "Bibliography and Resources" heading (after the
#     blank paragraph that follows the Double Buffering discussion).

$d = $word.ActiveDocument

function Get-ParagraphByPrefix($doc, [string]$prefix) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like ($prefix + "*")) {
            return $p
        }
    }
    return $null
}

function Get-ParagraphIndexByPrefix($doc, [string]$prefix) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like ($prefix + "*")) {
            return $idx
        }
    }
    return -1
}

# --- 1. "Benchmark" paragraph before "To benchmark and validate..." -------
$toPara = Get-ParagraphByPrefix $d "To benchmark and validate"
$toPara.Range.InsertParagraphBefore()
$toIdx = Get-ParagraphIndexByPrefix $d "To benchmark and validate"
$benchmarkPara = $d.Paragraphs.Item($toIdx - 1)
$benchmarkPara.Range.Text = "Benchmark"
$benchmarkPara.Range.Font.Underline = 1

# --- 2. "Validation" paragraph before "Furthermore, t..." ------------------
$furtherPara = Get-ParagraphByPrefix $d "Furthermore, t"
$furtherPara.Range.InsertParagraphBefore()
$furtherIdx = Get-ParagraphIndexByPrefix $d "Furthermore, t"
$validationPara = $d.Paragraphs.Item($furtherIdx - 1)
$validationPara.Range.Text = "Validation"
$validationPara.Range.Font.Underline = 1

# --- 3. "Results" heading before "Bibliography and Resources" --------------
$biblioPara = Get-ParagraphByPrefix $d "Bibliography and Resources"
$biblioPara.Range.InsertParagraphBefore()
$biblioIdx = Get-ParagraphIndexByPrefix $d "Bibliography and Resources"
$resultsPara = $d.Paragraphs.Item($biblioIdx - 1)
$resultsPara.Range.Text = "Results"

Write-Output "Done."
